# Applies the "Atualizacoes 16 de janeiro de 2024" updates to the
# ValueSet-DischDestVS workbook.

$wb = $excel.ActiveWorkbook

# --- Rename the second sheet ---------------------------------------------
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include from Discharge Locati"

# --- Update values on the Metadata sheet ----------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.1 -> 0.0.0
$wsMeta.Range("B3").Value = "0.0.0"

# Title: ValueSet of Discharge Location -> Discharge Location
$wsMeta.Range("B5").Value = "Discharge Location"

# Experimental value was blank, now set to the text "false".
# A bare .Value = "false" gets auto-coerced into a Boolean by the engine
# (the same way typing FALSE into a General formatted cell in real Excel
# does), so we write it through a quote-prefixed scratch cell and paste
# only the resulting text value into the target, which keeps the target
# cell's existing style/format untouched.
$scratch = $wsMeta.Cells.Item(100, 100)
$scratch.Value = "'false"
$scratch.Copy() | Out-Null
$wsMeta.Range("B7").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false

# Date: 2023-11-21T19:08:35-03:00 -> 2024-01-11T13:00:00-03:00
$wsMeta.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# Description updated wording
$wsMeta.Range("B12").Value = "ValueSet about discharge destination following acute care hospitalization. What type of place was the patient discharged to?"

# --- Update values on the Include sheet -----------------------------------
$wsInclude.Range("B10").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/DischDestCS"
